# Weekly refresh of the "Pepino ensalada" sheet: a new record is inserted
# at the top of the data block (row 288), pushing the existing rows
# 288-315 down to 289-316 (dimension grows from A1:R315 to A1:R316).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 288 - this shifts rows
# 288..315 down to 289..316 and extends the used range automatically.
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with the new weekly record.
$ws.Cells.Item(288, 1).Value = 5
$ws.Cells.Item(288, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(288, 3).Value = 'Maule'
$ws.Cells.Item(288, 4).Value = 44578
$ws.Cells.Item(288, 5).Value = 7
$ws.Cells.Item(288, 6).Value = 100112043
$ws.Cells.Item(288, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(288, 8).Value = 'Sin especificar'
$ws.Cells.Item(288, 9).Value = 'Primera'
$ws.Cells.Item(288, 10).Value = 500
$ws.Cells.Item(288, 11).Value = 8000
$ws.Cells.Item(288, 12).Value = 8000
$ws.Cells.Item(288, 13).Value = 8000
$ws.Cells.Item(288, 14).Value = '$/caja 80 unidades'
$ws.Cells.Item(288, 15).Value = 'Región del Maule'
$ws.Cells.Item(288, 16).Value = 100
$ws.Cells.Item(288, 17).Value = 80
$ws.Cells.Item(288, 18).Value = 'Hortaliza'
